$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 29166.5
$ws.Range("I21").Value = 36249.75
$ws.Range("J21").Value = 15000
$ws.Range("K21").Value = 36249.75
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = -35781.75
$ws.Range("N21").Value = -15936

$ws.Range("H23").Value = 29166.5
$ws.Range("I23").Value = 36249.75
$ws.Range("J23").Value = 15000
$ws.Range("K23").Value = 36249.75
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = -36015.75
$ws.Range("N23").Value = -15468

$ws.Range("H41").Value = 759
$ws.Range("I41").Value = 256.07144
$ws.Range("J41").Value = 2167.2
$ws.Range("K41").Value = 256.07144
$ws.Range("L41").Value = 2167.2
$ws.Range("M41").Value = 183.92856
$ws.Range("N41").Value = -3047.2

$ws.Range("H62").Value = 3405.375
$ws.Range("I62").Value = 2852.25
$ws.Range("J62").Value = 3958.5
$ws.Range("K62").Value = 2852.25
$ws.Range("L62").Value = 3958.5
$ws.Range("M62").Value = -2228.25
$ws.Range("N62").Value = -5206.5

$ws.Range("H65").Value = 3405.375
$ws.Range("I65").Value = 2852.25
$ws.Range("J65").Value = 3958.5
$ws.Range("K65").Value = 14261.25
$ws.Range("L65").Value = 19792.5
$ws.Range("M65").Value = -11141.25
$ws.Range("N65").Value = -26032.5

$ws.Range("H138").Value = 2088.9468
$ws.Range("I138").Value = 1149.7333
$ws.Range("J138").Value = 3497.7666
$ws.Range("K138").Value = 3449.199900000001
$ws.Range("L138").Value = 10493.2998
$ws.Range("M138").Value = 1690.800099999999
$ws.Range("N138").Value = -20773.2998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2040.5834
$ws.Range("I61").Value = 1788.1052
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1788.1052
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1576.1052
$ws.Range("N61").Value = -3424

$ws.Range("H136").Value = 2040.5834
$ws.Range("I136").Value = 1788.1052
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 5364.3156
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -2814.3156
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1824.1428
$ws.Range("I99").Value = 1214.75
$ws.Range("J99").Value = 2636.6667
$ws.Range("K99").Value = 1214.75
$ws.Range("L99").Value = 2636.6667
$ws.Range("M99").Value = 283.25
$ws.Range("N99").Value = -5632.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 191.25
$ws.Range("I19").Value = 191.25
$ws.Range("K19").Value = 191.25
$ws.Range("M19").Value = -21.25

$ws.Range("H22").Value = 397
$ws.Range("I22").Value = 249.625
$ws.Range("J22").Value = 691.75
$ws.Range("K22").Value = 249.625
$ws.Range("L22").Value = 691.75
$ws.Range("M22").Value = 100.375
$ws.Range("N22").Value = -1391.75

$ws.Range("H24").Value = 191.25
$ws.Range("I24").Value = 191.25
$ws.Range("K24").Value = 191.25
$ws.Range("M24").Value = -21.25

$ws.Range("H41").Value = 1439.75
$ws.Range("I41").Value = 1439.75
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1439.75
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()

$ws.Range("H50").Value = 29450
$ws.Range("J50").Value = 29450
$ws.Range("L50").Value = 29450
$ws.Range("N50").Value = -30700

$ws.Range("H51").Value = 11099
$ws.Range("J51").Value = 11099
$ws.Range("L51").Value = 11099
$ws.Range("N51").Value = -12571

$ws.Range("H59").Value = 21144.25
$ws.Range("I59").Value = 3000
$ws.Range("J59").Value = 39288.5
$ws.Range("K59").Value = 3000
$ws.Range("L59").Value = 39288.5
$ws.Range("M59").Value = -1855
$ws.Range("N59").Value = -41578.5

$ws.Range("H60").Value = 12782.4
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 13478
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 13478
$ws.Range("M60").Value = -9489
$ws.Range("N60").Value = -14500

$ws.Range("H61").Value = 11099
$ws.Range("J61").Value = 11099
$ws.Range("L61").Value = 11099
$ws.Range("N61").Value = -11795

$ws.Range("H68").Value = 24977
$ws.Range("J68").Value = 24977
$ws.Range("L68").Value = 24977
$ws.Range("N68").Value = -26475

$ws.Range("H71").Value = 24977
$ws.Range("J71").Value = 24977
$ws.Range("L71").Value = 74931
$ws.Range("N71").Value = -82419

$ws.Range("H74").Value = 33333.332
$ws.Range("J74").Value = 33333.332
$ws.Range("L74").Value = 33333.332
$ws.Range("N74").Value = -35081.332

$ws.Range("H77").Value = 33333.332
$ws.Range("J77").Value = 33333.332
$ws.Range("L77").Value = 99999.99600000001
$ws.Range("N77").Value = -108735.996

$ws.Range("H122").Value = 2048.9285
$ws.Range("I122").Value = 2039.5625
$ws.Range("J122").Value = 2061.4167
$ws.Range("K122").Value = 6118.6875
$ws.Range("L122").Value = 6184.250100000001
$ws.Range("M122").Value = -3668.6875
$ws.Range("N122").Value = -11084.2501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 772989.5
$ws.Range("I12").Value = 77.42856999999999
$ws.Range("J12").Value = 1073566.5
$ws.Range("K12").Value = 232.28571
$ws.Range("L12").Value = 3220699.5
$ws.Range("M12").Value = -59.28570999999999
$ws.Range("N12").Value = -3221045.5

$ws.Range("H134").Value = 3478.25
$ws.Range("I134").Value = 1829.9546
$ws.Range("J134").Value = 6068.4287
$ws.Range("K134").Value = 5489.8638
$ws.Range("L134").Value = 18205.2861
$ws.Range("M134").Value = -419.8638000000001
$ws.Range("N134").Value = -28345.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3749.9092
$ws.Range("I80").Value = 4167.5
$ws.Range("J80").Value = 3248.8
$ws.Range("K80").Value = 4167.5
$ws.Range("L80").Value = 3248.8
$ws.Range("M80").Value = -3169.5
$ws.Range("N80").Value = -5244.8

$ws.Range("H83").Value = 3749.9092
$ws.Range("I83").Value = 4167.5
$ws.Range("J83").Value = 3248.8
$ws.Range("K83").Value = 20837.5
$ws.Range("L83").Value = 16244
$ws.Range("M83").Value = -15845.5
$ws.Range("N83").Value = -26228

$ws.Range("H97").Value = 79073.75
$ws.Range("I97").Value = 89798.57000000001
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 89798.57000000001
$ws.Range("L97").Value = 4000
$ws.Range("M97").Value = -89302.57000000001
$ws.Range("N97").Value = -4992

$ws.Range("H132").Value = 1793.82
$ws.Range("I132").Value = 1218.8379
$ws.Range("J132").Value = 3430.3076
$ws.Range("K132").Value = 3656.5137
$ws.Range("L132").Value = 10290.9228
$ws.Range("M132").Value = -1126.5137
$ws.Range("N132").Value = -15350.9228
